$d = $word.ActiveDocument

# --- Edit 1: merge the split "Date de réception" runs into one ------------
# The original paragraph had 3 separate runs whose texts were
# "Date de réceptio" + "n" + " : {{object.date_reception...}}".
# Replace the concatenation of the first two ("Date de réceptio" + "n")
# with the correctly spelled word so the run collapses into a single run.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Date de réceptio" + "n" + " : ", $false, $false, $false, $false, $false, $true, 1, $false, "Date de réception : ", 2) | Out-Null

# --- Edit 2: drop the "Établissement (scénario retenu)" line --------------
# Keep the (now empty) paragraph but remove its text, and reset the
# paragraph's "space before" to 0 (it was 120).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tablissement (sc*nario retenu) : *conclusion_etablissement*") {
        $r = $p.Range
        $r.Text = ""
        $p.SpaceBefore = 0
        break
    }
}

# --- Edit 3: add a trailing space to the "Aliment" line and drop "Analyse"-
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Aliment (sc*nario retenu)*conclusion_aliment*") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = "Aliment (scénario retenu) : {{ object.conclusion_aliment or '-' }} "
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Analyse (sc*nario retenu)*conclusion_analyse*") {
        $p.Range.Delete()
        break
    }
}
